$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-categorize existing entries ---

# Row 46: SonarQube -> ring Adopt -> Trial
$ws.Range("B46").Value = "Trial"

# Row 47: Maestro -> ring Assess -> Trial
$ws.Range("B47").Value = "Trial"

# Row 48: OpenRewrite -> ring Assess -> Trial
$ws.Range("B48").Value = "Trial"

# Row 60: Cloud CI/CD -> ring Adopt -> Assess, add a description
$ws.Range("B60").Value = "Assess"
$ws.Range("E60").Value = "requirements not fulfilled for our team, otherwise a clear adopt"

# --- Add a new entry: Carthage (row 75) ---
$ws.Range("A75").Value = "Carthage"
$ws.Range("B75").Value = "Hold"
$ws.Range("C75").Value = "Tools"
$ws.Range("D75").Value = "FALSE"

# --- Update view state to match where the edits were made ---
$ws.Application.ActiveWindow.ScrollRow = 38
$ws.Range("E73").Select()
